# Lesson 3 Task 1
# Update the title on the first (title) slide of "03 CSS Basics.pptx"
# from "Client Side Scripting" to "JavaScript Programming".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the (center) title placeholder robustly, falling back to shape 1.
$title = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "Client Side Scripting") {
        $title = $sh
    }
}
if ($title -eq $null) {
    $title = $s.Shapes.Item(1)
}

$title.TextFrame.TextRange.Text = "JavaScript Programming"
